# previsao_retorno.xlsx - atualizei dados da bibi e add
# Refresh of the "Resumo_por_Cliente" report: the "situacao" column
# (INATIVO - X.X meses sem comprar) was recomputed against a newer
# reference date, nudging a batch of rows by +0.1 month, and two
# client rows (id 9247 / row 69, id 28458 / row 116) picked up fresh
# purchase-window metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- situacao (column J) recompute: "INATIVO - X.X meses sem comprar" ---
$ws.Cells.Item(5, 10).Value   = "INATIVO - 15.5 meses sem comprar"
$ws.Cells.Item(6, 10).Value   = "INATIVO - 17.0 meses sem comprar"
$ws.Cells.Item(16, 10).Value  = "INATIVO - 40.5 meses sem comprar"
$ws.Cells.Item(24, 10).Value  = "INATIVO - 38.5 meses sem comprar"
$ws.Cells.Item(25, 10).Value  = "INATIVO - 0.5 meses sem comprar"
$ws.Cells.Item(31, 10).Value  = "INATIVO - 7.5 meses sem comprar"
$ws.Cells.Item(39, 10).Value  = "INATIVO - 32.9 meses sem comprar"
$ws.Cells.Item(46, 10).Value  = "INATIVO - 6.8 meses sem comprar"
$ws.Cells.Item(47, 10).Value  = "INATIVO - 16.6 meses sem comprar"
$ws.Cells.Item(51, 10).Value  = "INATIVO - 8.2 meses sem comprar"
$ws.Cells.Item(67, 10).Value  = "INATIVO - 28.5 meses sem comprar"
$ws.Cells.Item(71, 10).Value  = "INATIVO - 12.0 meses sem comprar"
$ws.Cells.Item(78, 10).Value  = "INATIVO - 8.4 meses sem comprar"
$ws.Cells.Item(81, 10).Value  = "INATIVO - 6.8 meses sem comprar"
$ws.Cells.Item(90, 10).Value  = "INATIVO - 15.7 meses sem comprar"
$ws.Cells.Item(92, 10).Value  = "INATIVO - 12.4 meses sem comprar"
$ws.Cells.Item(93, 10).Value  = "INATIVO - 11.8 meses sem comprar"
$ws.Cells.Item(100, 10).Value = "INATIVO - 33.6 meses sem comprar"
$ws.Cells.Item(104, 10).Value = "INATIVO - 37.9 meses sem comprar"
$ws.Cells.Item(106, 10).Value = "INATIVO - 15.2 meses sem comprar"

# --- row 69 (id_cliente 9247, DAIANA HELENA PEREIRA DA SILVA): refreshed window ---
$ws.Cells.Item(69, 2).Value = 0.42
$ws.Cells.Item(69, 4).Value = 0.5
$ws.Cells.Item(69, 5).Value = 34
$ws.Cells.Item(69, 6).Value = 0.5
$ws.Cells.Item(69, 8).Value = 45855.7533912037
$ws.Cells.Item(69, 9).Value = 45886.7533912037

# --- row 116 (id_cliente 28458, BEMOL S/A): refreshed window ---
$ws.Cells.Item(116, 5).Value = 16733
$ws.Cells.Item(116, 8).Value = 45855.65665509259
$ws.Cells.Item(116, 9).Value = 45856.65665509259
